$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Python/IDLE version string (row 4, column B)
$ws.Range("B4").Value = "3.8.10 (Initial development on 3.8.3rc1)"

# Widen column B
$ws.Columns("B").ColumnWidth = 34.5546875

# Move the active selection from D12 to E10
$ws.Range("E10").Select()
